$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.332.59'
$ws.Range("E2").Value = '  -0.77%  '
$ws.Range("D3").Value = '1.709.96'
$ws.Range("E3").Value = '  -0.97%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.ClearFormats()
$ws.Range("E4").Value = '  -0.11%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '224.35'
$c.ClearFormats()
$ws.Range("E5").Value = '  -0.63%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.5297'
$c.ClearFormats()
$ws.Range("E6").Value = '  -1.23%  '
$ws.Range("E7").Value = '  -0.06%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.2665'
$c.ClearFormats()
$ws.Range("E8").Value = '  -0.19%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.06623'
$c.ClearFormats()
$ws.Range("E9").Value = '  +0.29%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '20.85'
$c.ClearFormats()
$ws.Range("E10").Value = '  -4.23%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07675'
$c.ClearFormats()
$ws.Range("E11").Value = '  -0.48%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '4.511'
$c.ClearFormats()
$ws.Range("E12").Value = '  -2.19%  '
$ws.Range("D13").Value = '1.943.34'
$ws.Range("E13").Value = '  -1.03%  '
$ws.Range("D14").Value = '1.708.43'
$ws.Range("E14").Value = '  -1.05%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.5816'
$c.ClearFormats()
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("D16").Value = '0.0₅8179'
$ws.Range("E16").Value = '  -1.44%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '67.83'
$c.ClearFormats()
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("D18").Value = '27.338.22'
$ws.Range("E18").Value = '  -0.82%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '217.76'
$c.ClearFormats()
$ws.Range("E19").Value = '  -1.21%  '
$ws.Range("E20").Value = '  -0.08%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '4.628'
$c.ClearFormats()
$ws.Range("E21").Value = '  -2.23%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '10.43'
$c.ClearFormats()
$ws.Range("E22").Value = '  -2.13%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '5.980'
$c.ClearFormats()
$ws.Range("E23").Value = '  -1.90%  '
$ws.Range("E24").Value = '  -0.10%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '143.82'
$c.ClearFormats()
$ws.Range("E25").Value = '  -3.06%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '1.687'
$c.ClearFormats()
$ws.Range("E26").Value = '  -1.26%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.1207'
$c.ClearFormats()
$ws.Range("E27").Value = '  -2.35%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '7.250'
$c.ClearFormats()
$ws.Range("E28").Value = '  -2.19%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '16.25'
$c.ClearFormats()
$ws.Range("E29").Value = '  -2.61%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.05364'
$c.ClearFormats()
$ws.Range("E30").Value = '  -3.69%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.294'
$c.ClearFormats()
$ws.Range("E31").Value = '  -0.67%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.465'
$c.ClearFormats()
$ws.Range("E32").Value = '  -2.48%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '3.424'
$c.ClearFormats()
$ws.Range("E33").Value = '  -1.06%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.648'
$c.ClearFormats()
$ws.Range("E34").Value = '  -0.77%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '2.863'
$c.ClearFormats()
$ws.Range("E35").Value = '  +1.32%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.9520'
$c.ClearFormats()
$ws.Range("E36").Value = '  -1.07%  '
$ws.Range("E37").Value = '  -1.31%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.5867'
$c.ClearFormats()
$ws.Range("E38").Value = '  -1.40%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.01639'
$c.ClearFormats()
$ws.Range("D40").Value = '1.068.40'
$ws.Range("E40").Value = '  +1.00%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '5.810'
$c.ClearFormats()
$ws.Range("E41").Value = '  -2.06%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.8441'
$c.ClearFormats()
$ws.Range("E42").Value = '  -1.45%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.004'
$c.ClearFormats()
$ws.Range("E43").Value = '  -0.03%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '101.07'
$c.ClearFormats()
$ws.Range("E44").Value = '  -0.39%  '
$ws.Range("D45").Value = '1.852.16'
$ws.Range("E45").Value = '  -0.96%  '
$ws.Range("D46").Value = '0.0₈117'
$ws.Range("E46").Value = '  +2.56%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '57.81'
$c.ClearFormats()
$ws.Range("E47").Value = '  -2.15%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.4523'
$c.ClearFormats()
$ws.Range("E48").Value = '  +1.91%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.004'
$c.ClearFormats()
$ws.Range("E49").Value = '  +0.27%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '8.075'
$c.ClearFormats()
$ws.Range("E50").Value = '  -1.67%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.05230'
$c.ClearFormats()
$ws.Range("E51").Value = '  -0.72%  '
